# Add functionality for ontogeny definition in a population.
# Adds two new header columns ("Ontogeny" and "Protein") to the
# Demographics sheet of the Population Parameters workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# New headers. Column R (18) gets "Ontogeny", column Q (17) gets "Protein".
# Written in this order so the new shared-string table entries land as
# 27 = "Ontogeny" and 28 = "Protein".
$ws.Cells.Item(1, 18).Value = "Ontogeny"
$ws.Cells.Item(1, 17).Value = "Protein"

# Move the active selection to Q2, matching the saved selection state.
$null = $ws.Range("Q2").Select()
